# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-job sheets with newly fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1267.58
$ws.Range("J17").Value = 1252.6327
$ws.Range("L17").Value = 3757.8981
$ws.Range("N17").Value = -4093.8981

$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20688

$ws.Range("H32").Value = 33336154
$ws.Range("J32").Value = 12503173
$ws.Range("L32").Value = 12503173
$ws.Range("N32").Value = -12503825

$ws.Range("H40").Value = 3449.6667
$ws.Range("I40").Value = 3435.4285
$ws.Range("J40").Value = 3499.5
$ws.Range("K40").Value = 3435.4285
$ws.Range("L40").Value = 3499.5
$ws.Range("M40").Value = -3260.4285
$ws.Range("N40").Value = -3849.5

$ws.Range("H63").Value = 59999
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 59999
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H76").Value = 4039
$ws.Range("I76").Value = 3897.8572
$ws.Range("J76").Value = 4088.4
$ws.Range("K76").Value = 3897.8572
$ws.Range("L76").Value = 4088.4
$ws.Range("M76").Value = -3582.8572
$ws.Range("N76").Value = -4718.4

$ws.Range("H79").Value = 4039
$ws.Range("I79").Value = 3897.8572
$ws.Range("J79").Value = 4088.4
$ws.Range("K79").Value = 3897.8572
$ws.Range("L79").Value = 4088.4
$ws.Range("M79").Value = -2805.8572
$ws.Range("N79").Value = -6272.4

$ws.Range("H86").Value = 285722140
$ws.Range("I86").Value = 666669000
$ws.Range("K86").Value = 666669000
$ws.Range("M86").Value = -666667877

$ws.Range("H89").Value = 285722140
$ws.Range("I89").Value = 666669000
$ws.Range("K89").Value = 3333345000
$ws.Range("M89").Value = -3333339384

$ws.Range("H107").Value = 16133282
$ws.Range("I107").Value = 18522432
$ws.Range("J107").Value = 6520.25
$ws.Range("K107").Value = 18522432
$ws.Range("L107").Value = 6520.25
$ws.Range("M107").Value = -18520512
$ws.Range("N107").Value = -10360.25

$ws.Range("H137").Value = 3695.868
$ws.Range("I137").Value = 4393.684
$ws.Range("K137").Value = 13181.052
$ws.Range("M137").Value = -10631.052

$ws.Range("H139").Value = 90265
$ws.Range("J139").Value = 90265
$ws.Range("L139").Value = 90265
$ws.Range("N139").Value = -100545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 59899
$ws.Range("J7").Value = 59899
$ws.Range("L7").Value = 59899
$ws.Range("N7").Value = -60127

$ws.Range("H61").Value = 4446.1113
$ws.Range("I61").Value = 4491.9116
$ws.Range("J61").Value = 3667.5
$ws.Range("K61").Value = 4491.9116
$ws.Range("L61").Value = 3667.5
$ws.Range("M61").Value = -4279.9116
$ws.Range("N61").Value = -4091.5

$ws.Range("H63").Value = 3606.8572
$ws.Range("I63").Value = 2166.3333
$ws.Range("J63").Value = 3999.7273
$ws.Range("K63").Value = 2166.3333
$ws.Range("L63").Value = 3999.7273
$ws.Range("M63").Value = -1480.3333
$ws.Range("N63").Value = -5371.7273

$ws.Range("H66").Value = 3606.8572
$ws.Range("I66").Value = 2166.3333
$ws.Range("J66").Value = 3999.7273
$ws.Range("K66").Value = 10831.6665
$ws.Range("L66").Value = 19998.6365
$ws.Range("M66").Value = -7399.666499999999
$ws.Range("N66").Value = -26862.6365

$ws.Range("H122").Value = 16455.174
$ws.Range("J122").Value = 43862
$ws.Range("L122").Value = 131586
$ws.Range("N122").Value = -136486

$ws.Range("H132").Value = 2189.2727
$ws.Range("I132").Value = 2054.889
$ws.Range("J132").Value = 2794
$ws.Range("K132").Value = 6164.667
$ws.Range("L132").Value = 8382
$ws.Range("M132").Value = -3634.667
$ws.Range("N132").Value = -13442

$ws.Range("H136").Value = 4446.1113
$ws.Range("I136").Value = 4491.9116
$ws.Range("J136").Value = 3667.5
$ws.Range("K136").Value = 13475.7348
$ws.Range("L136").Value = 11002.5
$ws.Range("M136").Value = -10925.7348
$ws.Range("N136").Value = -16102.5

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6125.0586
$ws.Range("I20").Value = 5108
$ws.Range("K20").Value = 5108
$ws.Range("M20").Value = -4861

$ws.Range("H82").Value = 11392.143
$ws.Range("I82").Value = 4957.6665
$ws.Range("K82").Value = 4957.6665
$ws.Range("M82").Value = -4574.6665

$ws.Range("H85").Value = 11392.143
$ws.Range("I85").Value = 4957.6665
$ws.Range("K85").Value = 4957.6665
$ws.Range("M85").Value = -3631.6665

$ws.Range("H115").Value = 60000
$ws.Range("J115").Value = 60000
$ws.Range("L115").Value = 60000
$ws.Range("N115").Value = -63134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1900.75
$ws.Range("J16").Value = 2292.889
$ws.Range("L16").Value = 2292.889
$ws.Range("N16").Value = -2866.889

$ws.Range("H113").Value = 1900.75
$ws.Range("J113").Value = 2292.889
$ws.Range("L113").Value = 2292.889
$ws.Range("N113").Value = -6632.889

$ws.Range("H120").Value = 29899
$ws.Range("J120").Value = 29899
$ws.Range("L120").Value = 29899
$ws.Range("N120").Value = -37157

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3667.3333
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 44991
$ws.Range("N132").Value = -50051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6482.5
$ws.Range("I102").Value = 6311.467
$ws.Range("K102").Value = 6311.467
$ws.Range("M102").Value = -4689.467

$ws.Range("H122").Value = 2549.75
$ws.Range("J122").Value = 3589.5
$ws.Range("L122").Value = 10768.5
$ws.Range("N122").Value = -15668.5

$ws.Range("H132").Value = 4291.3335
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4291.3335
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12874.0005
$ws.Range("N132").Value = -17934.0005
$ws.Range("M132").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H137").Value = 40709
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H140").Value = 107474.5
$ws.Range("J140").Value = 107474.5
$ws.Range("L140").Value = 107474.5
$ws.Range("N140").Value = -117834.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4242.2
$ws.Range("I40").Value = 4163.8125
$ws.Range("J40").Value = 4555.75
$ws.Range("K40").Value = 4163.8125
$ws.Range("L40").Value = 4555.75
$ws.Range("M40").Value = -4027.8125
$ws.Range("N40").Value = -4827.75

$ws.Range("H122").Value = 5533.75
$ws.Range("I122").Value = 5835.8823
$ws.Range("K122").Value = 17507.6469
$ws.Range("M122").Value = -15057.6469

$ws.Range("H132").Value = 3918.4119
$ws.Range("I132").Value = 3918.4119
$ws.Range("K132").Value = 11755.2357
$ws.Range("M132").Value = -9225.235700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 41952.43
$ws.Range("J46").Value = 41952.43
$ws.Range("L46").Value = 41952.43
$ws.Range("N46").Value = -42414.43

$ws.Range("H64").Value = 59999
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59999
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59999
$ws.Range("N64").Value = -60495
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 59999
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59999
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59999
$ws.Range("N67").Value = -61715
$ws.Range("M67").ClearContents()

$ws.Range("H134").Value = 41952.43
$ws.Range("J134").Value = 41952.43
$ws.Range("L134").Value = 125857.29
$ws.Range("N134").Value = -130927.29
